# "modified code to use port expander for most digital signals"
# Replace the direct Arduino digital-pin references (D0-D13, A3-A5) used for
# most of the relay-driver signals with the MCP23017 port-expander pins
# (GPA0-GPA7, GPB0-GPB7), keeping only D2/D3 (button / tune button) and
# VCC/GND as direct Arduino pins.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: swap Dn/An pin names for GPAn/GPBn port-expander pins -------
$ws.Range("A2").Value  = "GPA0"
$ws.Range("A3").Value  = "GPA1"
$ws.Range("A4").Value  = "GPA2"
$ws.Range("A5").Value  = "GPA3"
$ws.Range("A6").Value  = "GPA4"
$ws.Range("A7").Value  = "GPA5"
$ws.Range("A8").Value  = "GPA6"
$ws.Range("A10").Value = "GPA7"
$ws.Range("A11").Value = "GPB0"
$ws.Range("A12").Value = "GPB1"
$ws.Range("A13").Value = "GPB2"
$ws.Range("A14").Value = "GPB3"
$ws.Range("A15").Value = "GPB4"
$ws.Range("A16").Value = "GPB5"
$ws.Range("A18").Value = "GPB6"
$ws.Range("A19").Value = "GPB7"

# A11 ("GPB0") keeps the "B0" portion bold, as a rich-text run.
$ws.Range("A11").Characters(3, 2).Font.Bold = $true

# A12 ("GPB1") is rendered fully bold (its whole cell font is bold).
$ws.Range("A12").Font.Bold = $true

# --- Remove the old "A5 / LED" row (row 20) and renumber the pin used -------
$ws.Rows(20).Delete()

# Make room for the LED row (now referencing D2) and a new "Tune button" row
# (D3) between the GPIO table and the VCC/GND rows.
$ws.Rows(21).Insert()
$ws.Rows(22).Insert()
$ws.Rows(23).Insert()

$ws.Range("A21").Value = "D2"
$ws.Range("B21").Value = 4
$ws.Range("C21").Value = "LED"
$ws.Range("F21").Value = 17

$ws.Range("A22").Value = "D3"
$ws.Range("E22").Value = "Tune button"

# The value wired to K9 changed description.
$ws.Range("E18").Value = "Cin/Cout"

# --- View: select A18 as the active cell -----------------------------------
$ws.Range("A18").Select()
